$wb = $excel.ActiveWorkbook

# ALC sheet
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1809.5294
$ws.Range("I92").Value = 2027.091
$ws.Range("J92").Value = 1410.6666
$ws.Range("K92").Value = 2027.091
$ws.Range("L92").Value = 1410.6666
$ws.Range("M92").Value = -779.0909999999999
$ws.Range("N92").Value = -3906.6666
$ws.Range("H112").Value = 1299.9811
$ws.Range("I112").Value = 737.375
$ws.Range("K112").Value = 2212.125
$ws.Range("M112").Value = -1104.125
$ws.Range("H116").Value = 509330.1
$ws.Range("I116").Value = 1114543.5
$ws.Range("J116").Value = 14155.546
$ws.Range("K116").Value = 1114543.5
$ws.Range("L116").Value = 14155.546
$ws.Range("M116").Value = -1111101.5
$ws.Range("N116").Value = -21039.546
$ws.Range("H129").Value = 915.03094
$ws.Range("I129").Value = 441.7143
$ws.Range("J129").Value = 951.8444
$ws.Range("K129").Value = 1325.1429
$ws.Range("L129").Value = 2855.5332
$ws.Range("M129").Value = 3674.8571
$ws.Range("N129").Value = -12855.5332
$ws.Range("H132").Value = 298929.88
$ws.Range("I132").Value = 5270.24
$ws.Range("J132").Value = 1114651.1
$ws.Range("K132").Value = 15810.72
$ws.Range("L132").Value = 3343953.3
$ws.Range("M132").Value = -13280.72
$ws.Range("N132").Value = -3349013.3
$ws.Range("H137").Value = 4037.5557
$ws.Range("I137").Value = 2191.5454
$ws.Range("J137").Value = 6938.4287
$ws.Range("K137").Value = 6574.6362
$ws.Range("L137").Value = 20815.2861
$ws.Range("M137").Value = -4024.6362
$ws.Range("N137").Value = -25915.2861
$ws.Range("H138").Value = 2387.03
$ws.Range("I138").Value = 805.6429000000001
$ws.Range("J138").Value = 3002.014
$ws.Range("K138").Value = 2416.9287
$ws.Range("L138").Value = 9006.042000000001
$ws.Range("M138").Value = 2723.0713
$ws.Range("N138").Value = -19286.042

# ARM sheet
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1026.2174
$ws.Range("I61").Value = 725.05554
$ws.Range("J61").Value = 2110.4
$ws.Range("K61").Value = 725.05554
$ws.Range("L61").Value = 2110.4
$ws.Range("M61").Value = -513.05554
$ws.Range("N61").Value = -2534.4
$ws.Range("H74").Value = 4180.2593
$ws.Range("I74").Value = 4584.4
$ws.Range("J74").Value = 3025.5715
$ws.Range("K74").Value = 4584.4
$ws.Range("L74").Value = 3025.5715
$ws.Range("M74").Value = -3710.4
$ws.Range("N74").Value = -4773.5715
$ws.Range("H77").Value = 4180.2593
$ws.Range("I77").Value = 4584.4
$ws.Range("J77").Value = 3025.5715
$ws.Range("K77").Value = 22922
$ws.Range("L77").Value = 15127.8575
$ws.Range("M77").Value = -18554
$ws.Range("N77").Value = -23863.8575
$ws.Range("H102").Value = 1836.875
$ws.Range("I102").Value = 1683.1666
$ws.Range("J102").Value = 2298
$ws.Range("K102").Value = 1683.1666
$ws.Range("L102").Value = 2298
$ws.Range("M102").Value = -61.16660000000002
$ws.Range("N102").Value = -5542
$ws.Range("H122").Value = 1976.1666
$ws.Range("I122").Value = 1295.5625
$ws.Range("J122").Value = 3337.375
$ws.Range("K122").Value = 3886.6875
$ws.Range("L122").Value = 10012.125
$ws.Range("M122").Value = -1436.6875
$ws.Range("N122").Value = -14912.125
$ws.Range("H132").Value = 2547
$ws.Range("I132").Value = 1183.8823
$ws.Range("J132").Value = 5121.778
$ws.Range("K132").Value = 3551.6469
$ws.Range("L132").Value = 15365.334
$ws.Range("M132").Value = -1021.6469
$ws.Range("N132").Value = -20425.334
$ws.Range("H136").Value = 1026.2174
$ws.Range("I136").Value = 725.05554
$ws.Range("J136").Value = 2110.4
$ws.Range("K136").Value = 2175.16662
$ws.Range("L136").Value = 6331.200000000001
$ws.Range("M136").Value = 374.83338
$ws.Range("N136").Value = -11431.2

# BSM sheet
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1572.4131
$ws.Range("I134").Value = 1015.97296
$ws.Range("J134").Value = 3860
$ws.Range("K134").Value = 3047.91888
$ws.Range("L134").Value = 11580
$ws.Range("M134").Value = -512.9188799999997
$ws.Range("N134").Value = -16650

# CRP sheet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38465850
$ws.Range("I31").Value = 1222.2
$ws.Range("J31").Value = 62506240
$ws.Range("K31").Value = 1222.2
$ws.Range("L31").Value = 62506240
$ws.Range("M31").Value = -927.2
$ws.Range("N31").Value = -62506830
$ws.Range("H34").Value = 38465850
$ws.Range("I34").Value = 1222.2
$ws.Range("J34").Value = 62506240
$ws.Range("K34").Value = 1222.2
$ws.Range("L34").Value = 62506240
$ws.Range("M34").Value = -1020.2
$ws.Range("N34").Value = -62506644
$ws.Range("H58").Value = 1790.1559
$ws.Range("I58").Value = 1545.5074
$ws.Range("J58").Value = 3429.3
$ws.Range("K58").Value = 1545.5074
$ws.Range("L58").Value = 3429.3
$ws.Range("M58").Value = -1342.5074
$ws.Range("N58").Value = -3835.3
$ws.Range("H123").Value = 39835
$ws.Range("J123").Value = 39835
$ws.Range("L123").Value = 39835
$ws.Range("N123").Value = -49635
$ws.Range("H132").Value = 3818.75
$ws.Range("I132").Value = 3075.125
$ws.Range("K132").Value = 9225.375
$ws.Range("M132").Value = -6695.375
$ws.Range("H134").Value = 7278.579
$ws.Range("I134").Value = 8599.923000000001
$ws.Range("J134").Value = 4415.6665
$ws.Range("K134").Value = 25799.769
$ws.Range("L134").Value = 13246.9995
$ws.Range("M134").Value = -23264.769
$ws.Range("N134").Value = -18316.9995
$ws.Range("H136").Value = 1790.1559
$ws.Range("I136").Value = 1545.5074
$ws.Range("J136").Value = 3429.3
$ws.Range("K136").Value = 4636.522199999999
$ws.Range("L136").Value = 10287.9
$ws.Range("M136").Value = -2086.522199999999
$ws.Range("N136").Value = -15387.9

# CUL sheet
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 826.5833
$ws.Range("I113").Value = 656.1539
$ws.Range("K113").Value = 1968.4617
$ws.Range("M113").Value = 201.5382999999999
$ws.Range("H131").Value = 12821429
$ws.Range("J131").Value = 1009.32355
$ws.Range("L131").Value = 3027.97065
$ws.Range("N131").Value = -13107.97065
$ws.Range("H137").Value = 3828.5833
$ws.Range("J137").Value = 3994.818
$ws.Range("L137").Value = 11984.454
$ws.Range("N137").Value = -22184.454

# GSM sheet
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 632.5714
$ws.Range("I97").Value = 767.5
$ws.Range("J97").Value = 452.66666
$ws.Range("K97").Value = 767.5
$ws.Range("L97").Value = 452.66666
$ws.Range("M97").Value = -271.5
$ws.Range("N97").Value = -1444.66666
$ws.Range("H102").Value = 2617.6
$ws.Range("I102").Value = 1802.4
$ws.Range("J102").Value = 4248
$ws.Range("K102").Value = 1802.4
$ws.Range("L102").Value = 4248
$ws.Range("M102").Value = -180.4000000000001
$ws.Range("N102").Value = -7492
$ws.Range("H132").Value = 3526.7144
$ws.Range("I132").Value = 2054.0833
$ws.Range("J132").Value = 5490.222
$ws.Range("K132").Value = 6162.249899999999
$ws.Range("L132").Value = 16470.666
$ws.Range("M132").Value = -3632.249899999999
$ws.Range("N132").Value = -21530.666

# LTW sheet
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 41754.445
$ws.Range("J128").Value = 41754.445
$ws.Range("L128").Value = 41754.445
$ws.Range("N128").Value = -51714.445
$ws.Range("H132").Value = 3262.1875
$ws.Range("I132").Value = 1379.6285
$ws.Range("J132").Value = 8330.615
$ws.Range("K132").Value = 4138.8855
$ws.Range("L132").Value = 24991.845
$ws.Range("M132").Value = -1608.8855
$ws.Range("N132").Value = -30051.845

# WVR sheet
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 144371920
$ws.Range("I96").Value = 252650000
$ws.Range("J96").Value = 1156.3334
$ws.Range("K96").Value = 252650000
$ws.Range("L96").Value = 1156.3334
$ws.Range("M96").Value = 252650000
$ws.Range("N96").Value = -3902.3334
